# Applies the "Initial Data File Updated" commit: appends 12 new
# transaction rows (60-71) to the "Transacciones" sheet, fixes a
# mis-categorised row (D56), and moves the selection/scroll position to
# the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

# ---------------------------------------------------------------------
# 1) Row 56 was re-categorised from "Despensa" to "Misc".
# ---------------------------------------------------------------------
$ws.Cells.Item(56, 4).Value = "Misc"

# ---------------------------------------------------------------------
# 2) New transaction rows 60-71.
# ---------------------------------------------------------------------

# -- Column A: dates (all written as the same serial date, then the
#    existing date NumberFormat from A59 is copied onto them so they
#    keep the "m/d/yyyy"-style formatting used by the rest of the
#    column instead of the engine inventing a brand-new style).
$ws.Cells.Item(60, 1).Value = 43573
$ws.Cells.Item(61, 1).Value = 43573
$ws.Cells.Item(62, 1).Value = 43573
$ws.Cells.Item(63, 1).Value = 43573
$ws.Cells.Item(64, 1).Value = 43573
$ws.Cells.Item(65, 1).Value = 43573
$ws.Cells.Item(66, 1).Value = 43573
$ws.Cells.Item(67, 1).Value = 43573
$ws.Cells.Item(68, 1).Value = 43574
$ws.Cells.Item(69, 1).Value = 43574
$ws.Cells.Item(70, 1).Value = 43574
$ws.Cells.Item(71, 1).Value = 43574
$ws.Range("A59").Copy()
$ws.Range("A60:A71").PasteSpecial(-4122)

# -- Column B: the literal peso amounts that drive each row's K/M math.
$ws.Cells.Item(60, 2).Value = 186
$ws.Cells.Item(61, 2).Value = 62
$ws.Cells.Item(62, 2).Value = 40
$ws.Cells.Item(63, 2).Value = 83
$ws.Cells.Item(64, 2).Value = 20
$ws.Cells.Item(65, 2).Value = 10
$ws.Cells.Item(66, 2).Value = 20
$ws.Cells.Item(67, 2).Value = 102
$ws.Cells.Item(68, 2).Value = 50
$ws.Cells.Item(69, 2).Value = 100
$ws.Cells.Item(70, 2).Value = 10
$ws.Cells.Item(71, 2).Value = 10

# -- Columns C-G: description / category / in-out / payment / place.
$ws.Cells.Item(60, 3).Value = "Boletos Cinemex"
$ws.Cells.Item(60, 4).Value = "Ocio"
$ws.Cells.Item(60, 5).Value = "Gasto"
$ws.Cells.Item(60, 6).Value = "Tarjeta Banamex"
$ws.Cells.Item(60, 7).Value = "Cinemex"

$ws.Cells.Item(61, 3).Value = "Boleto Cinemex"
$ws.Cells.Item(61, 4).Value = "Pagos"
$ws.Cells.Item(61, 5).Value = "Ingreso"
$ws.Cells.Item(61, 6).Value = "Efectivo"
$ws.Cells.Item(61, 7).Value = "NA"

$ws.Cells.Item(62, 3).Value = "Boleto Cinemex"
$ws.Cells.Item(62, 4).Value = "Pagos"
$ws.Cells.Item(62, 5).Value = "Ingreso"
$ws.Cells.Item(62, 6).Value = "Efectivo"
$ws.Cells.Item(62, 7).Value = "NA"

$ws.Cells.Item(63, 3).Value = "Comida China"
$ws.Cells.Item(63, 4).Value = "Comida"
$ws.Cells.Item(63, 5).Value = "Gasto"
$ws.Cells.Item(63, 6).Value = "Efectivo"
$ws.Cells.Item(63, 7).Value = "Via Alta"

$ws.Cells.Item(64, 3).Value = "Estacionamiento"
$ws.Cells.Item(64, 4).Value = "Estacionamiento"
$ws.Cells.Item(64, 5).Value = "Gasto"
$ws.Cells.Item(64, 6).Value = "Efectivo"
$ws.Cells.Item(64, 7).Value = "Via Alta"

$ws.Cells.Item(65, 3).Value = "Pan Mil Hojas"
$ws.Cells.Item(65, 4).Value = "Golosinas"
$ws.Cells.Item(65, 5).Value = "Gasto"
$ws.Cells.Item(65, 6).Value = "Efectivo"
$ws.Cells.Item(65, 7).Value = "Salamanca Centro"

$ws.Cells.Item(66, 3).Value = "Estacionamiento"
$ws.Cells.Item(66, 4).Value = "Estacionamiento"
$ws.Cells.Item(66, 5).Value = "Gasto"
$ws.Cells.Item(66, 6).Value = "Efectivo"
$ws.Cells.Item(66, 7).Value = "Salamanca Centro"

$ws.Cells.Item(67, 3).Value = "Tacos"
$ws.Cells.Item(67, 4).Value = "Comida"
$ws.Cells.Item(67, 5).Value = "Gasto"
$ws.Cells.Item(67, 6).Value = "Efectivo"
$ws.Cells.Item(67, 7).Value = "Salamanca Centro"

$ws.Cells.Item(68, 3).Value = "Gasolina - Chore"
$ws.Cells.Item(68, 4).Value = "Gasolina"
$ws.Cells.Item(68, 5).Value = "Ingreso"
$ws.Cells.Item(68, 6).Value = "Efectivo"
$ws.Cells.Item(68, 7).Value = "NA"

$ws.Cells.Item(69, 3).Value = "Gasolina - Hugo"
$ws.Cells.Item(69, 4).Value = "Gasolina"
$ws.Cells.Item(69, 5).Value = "Ingreso"
$ws.Cells.Item(69, 6).Value = "Efectivo"
$ws.Cells.Item(69, 7).Value = "NA"

$ws.Cells.Item(70, 3).Value = "Propina - Limpiaparabrisas"
$ws.Cells.Item(70, 4).Value = "Misc"
$ws.Cells.Item(70, 5).Value = "Gasto"
$ws.Cells.Item(70, 6).Value = "Efectivo"
$ws.Cells.Item(70, 7).Value = "Leon Centro"

$ws.Cells.Item(71, 3).Value = "Ingreso de monedas a alcancía"
$ws.Cells.Item(71, 4).Value = "ahorro"
$ws.Cells.Item(71, 5).Value = "Gasto"
$ws.Cells.Item(71, 6).Value = "Efectivo"
$ws.Cells.Item(71, 7).Value = "Alcancía"

# -- Column K ("Monto Actual" running cash balance): row 60 recomputes
#    from the previous balance, the rest repeat the same literal.
$ws.Range("K60").Formula = "=K59-B60"
$ws.Cells.Item(61, 11).Value = 7900.24
$ws.Cells.Item(62, 11).Value = 7900.24
$ws.Cells.Item(63, 11).Value = 7900.24
$ws.Cells.Item(64, 11).Value = 7900.24
$ws.Cells.Item(65, 11).Value = 7900.24
$ws.Cells.Item(66, 11).Value = 7900.24
$ws.Cells.Item(67, 11).Value = 7900.24
$ws.Cells.Item(68, 11).Value = 7900.24
$ws.Cells.Item(69, 11).Value = 7900.24
$ws.Cells.Item(70, 11).Value = 7900.24
$ws.Cells.Item(71, 11).Value = 7900.24

# -- Column L ("Santander" balance): unchanged literal for every row.
for ($r = 60; $r -le 71; $r++) {
    $ws.Cells.Item($r, 12).Value = 2527.5700000000002
}

# -- Column M (loose-change tally): row 60 is a literal carry-over,
#    every other row nudges the previous row by that row's B amount.
$ws.Cells.Item(60, 13).Value = 175
$ws.Range("M61").Formula = "=M60+B61"
$ws.Range("M62").Formula = "=M61+B62"
$ws.Range("M63").Formula = "=M62-B63"
$ws.Range("M64").Formula = "=M63-B64"
$ws.Range("M65").Formula = "=M64-B65"
$ws.Range("M66").Formula = "=M65-B66"
$ws.Range("M67").Formula = "=M66-B67"
$ws.Range("M68").Formula = "=M67+B68"
$ws.Range("M69").Formula = "=M68+B69"
$ws.Range("M70").Formula = "=M69-B70"
$ws.Range("M71").Formula = "=M70-B71"

# -- Columns N/O (period subtotal + running total minus the 4000
#    baseline). Written range-at-a-time so that rows which were
#    filled down together in the original workbook are re-serialised
#    as one OOXML shared formula, matching how Excel would have
#    produced them.
$ws.Range("N60").Formula = "=SUM(K60:M60)"
$ws.Range("O60").Formula = "=N60-4000"

$ws.Range("N61:N63").Formula = "=SUM(K61:M61)"
$ws.Range("O61:O63").Formula = "=N61-4000"

$ws.Range("N64").Formula = "=SUM(K64:M64)"
$ws.Range("O64").Formula = "=N64-4000"

$ws.Range("N65").Formula = "=SUM(K65:M65)"
$ws.Range("O65").Formula = "=N65-4000"

$ws.Range("N66:N69").Formula = "=SUM(K66:M66)"
$ws.Range("O66:O69").Formula = "=N66-4000"

$ws.Range("N70").Formula = "=SUM(K70:M70)"
$ws.Range("O70").Formula = "=N70-4000"

$ws.Range("N71").Formula = "=SUM(K71:M71)"
$ws.Range("O71").Formula = "=N71-4000"

# ---------------------------------------------------------------------
# 3) Sheet dimension grows to A1:V71 automatically as cells are
#    written; move the selection to the new last entry row so the
#    workbook re-opens scrolled to where the new rows were added.
# ---------------------------------------------------------------------
$ws.Range("L69").Select()

Write-Host "Applied Initial Data File Updated edit"
